# feat: add 2022-Q1 data
#
# Inserts a new "2022-Q1" sheet (fund-holding detail, same layout as the
# 2021-Qx sheets) right before the "总计" (totals) sheet, and updates the
# "总计" sheet with a new leading row summarising 2022-Q1 (12 holdings,
# 0.86 亿元), shifting the previously-existing rows down.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Re-create the sheet order/ids so the final workbook looks like:
#      2021-Q2 (id1) 2021-Q3 (id2) 2021-Q4 (id3) 2022-Q1 (id4) 总计 (id5)
#    The cleanest way to land on that exact id sequence with this host is to
#    drop the old "总计" sheet (freeing id 4) and then clone the "2021-Q4"
#    sheet twice - the clones pick up fresh, sequential ids (4, then 5) and
#    they also inherit the right sheetPr/pageMargins boilerplate that the
#    hand-authored quarter sheets use (2021-Q2 alone has different margins).
# ---------------------------------------------------------------------------

$wb.Worksheets.Item("总计").Delete() | Out-Null

$template = $wb.Worksheets.Item("2021-Q4")
$template.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$q1 = $wb.Worksheets.Item($wb.Worksheets.Count)
$q1.Name = "2022-Q1"

$template2 = $wb.Worksheets.Item("2021-Q4")
$template2.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$total = $wb.Worksheets.Item($wb.Worksheets.Count)
$total.Name = "总计"

# ---------------------------------------------------------------------------
# 2. Build the "2022-Q1" sheet content.
#    Wipe the cloned "2021-Q4" data but keep the inherited formatting (the
#    bold-bordered style used for the header row / index column), then trim
#    the sheet back down to the 12 rows of real data (A1:H13).
# ---------------------------------------------------------------------------

$q1.Range("A1:H50").ClearContents()
$q1.Range("A14:A50").EntireRow.Delete()

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Columns B, D, E, F, G hold numeric-looking text ("006648", "6.22", ...)
# that must stay text (and keep leading zeros in the fund codes), so mark
# them as Text before writing - otherwise the host auto-coerces them to
# numbers.
$q1.Range("B2:B13").NumberFormat = "@"
$q1.Range("D2:G13").NumberFormat = "@"

$q1Rows = @(
    @("006648", "汇安多因子混合A",       "6.22", "93.26", "3.84", "0.2388", 5),
    @("009381", "汇安核心资产混合A",       "4.50", "94.22", "3.60", "0.1620", 8),
    @("013867", "汇安优势企业精选混合A",   "4.29", "91.60", "3.49", "0.1497", 4),
    @("010558", "汇安鑫利优选混合A",       "2.21", "92.83", "4.32", "0.0955", 3),
    @("006649", "汇安多因子混合C",       "2.26", "93.26", "3.84", "0.0868", 5),
    @("003889", "汇安丰泽灵活配置混合A",   "1.36", "93.28", "3.84", "0.0522", 5),
    @("010559", "汇安鑫利优选混合C",       "0.73", "92.83", "4.32", "0.0315", 3),
    @("003890", "汇安丰泽灵活配置混合C",   "0.51", "93.28", "3.84", "0.0196", 5),
    @("007775", "汇安量化先锋混合A",       "0.38", "94.51", "3.16", "0.0120", 7),
    @("013868", "汇安优势企业精选混合C",   "0.22", "91.60", "3.49", "0.0077", 4),
    @("007776", "汇安量化先锋混合C",       "0.11", "94.51", "3.16", "0.0035", 7),
    @("009382", "汇安核心资产混合C",       "0.08", "94.22", "3.60", "0.0029", 8)
)

$r = 2
foreach ($row in $q1Rows) {
    $q1.Cells.Item($r, 1).Value = ($r - 2)
    $q1.Cells.Item($r, 2).Value = $row[0]
    $q1.Cells.Item($r, 3).Value = $row[1]
    $q1.Cells.Item($r, 4).Value = $row[2]
    $q1.Cells.Item($r, 5).Value = $row[3]
    $q1.Cells.Item($r, 6).Value = $row[4]
    $q1.Cells.Item($r, 7).Value = $row[5]
    $q1.Cells.Item($r, 8).Value = $row[6]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# 3. Build the "总计" sheet content: the same 日期/持有数量(只)/持有市值(亿元)
#    table as before, with a new 2022-Q1 row inserted at the top and the
#    index column renumbered.
# ---------------------------------------------------------------------------

$total.Range("A1:H50").ClearContents()
$total.Range("E1:H50").EntireColumn.Delete()
$total.Range("A6:A50").EntireRow.Delete()

$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

$totalRows = @(
    @("2022-Q1", 12, 0.86),
    @("2021-Q4", 22, 2.38),
    @("2021-Q3", 30, 6.28),
    @("2021-Q2", 14, 2.03)
)

$r = 2
foreach ($row in $totalRows) {
    $total.Cells.Item($r, 1).Value = ($r - 2)
    $total.Cells.Item($r, 2).Value = $row[0]
    $total.Cells.Item($r, 3).Value = $row[1]
    $total.Cells.Item($r, 4).Value = $row[2]
    $r = $r + 1
}

Write-Output "2022-Q1 sheet added; 总计 sheet updated"
